$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("tomato millet with powder", "barnyard millet: 0.75 cup,mtr tomato rice powder: 2.0 pack,olive oil: 4.0 tsp,raw chana dal: 1.0 tbsp,peanuts: 2.5 tbsp", 3, 323.3333333333333, 7.5, 11.66666666666667, 42.33333333333334),
    @("coffee with oat creamer", "planet oat creamer: 2.0 tbsp", 1, 50, 0, 2, 8),
    @("red apple pear orange smoothie no yogurt", "apple: 1.0,green pear: 1.0,orange: 1.0,pitted date: 8.0,unsweetened almond milk: 1.25 cup", 2, 253.75, 1.125, 1.25, 58.625),
    @("capsicum besan sabji", "green capsicum: 1.0,besan: 0.33 cup,olive oil: 2.0 tsp,onion: 0.75", 1, 263.48, 19.18, 9.98, 25.24)
)

$startRow = 108
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
